$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A7").Value = "save.p"
$ws.Range("B11").Value = "species, occ"
